$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 234, shifting existing rows 234..344 down to 235..345
$ws.Rows("234:234").Insert()

# Populate the newly inserted row 234 with its data
$ws.Range("A234").Value = 5
$ws.Range("B234").Value = "Macroferia Regional de Talca"
$ws.Range("C234").Value = "Maule"
$ws.Range("D234").Value = 45141
$ws.Range("E234").Value = 7
$ws.Range("F234").Value = 100112017
$ws.Range("G234").Value = "Apio"
$ws.Range("H234").Value = "Americana (o)"
$ws.Range("I234").Value = "Primera"
$ws.Range("J234").Value = 700
$ws.Range("K234").Value = 5000
$ws.Range("L234").Value = 5000
$ws.Range("M234").Value = 5000
$ws.Range("N234").Value = "`$/docena de matas"
$ws.Range("O234").Value = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P234").Value = 833
$ws.Range("Q234").Value = 6
$ws.Range("R234").Value = "Hortaliza"
